# Add a new "GatewaySource" worksheet, positioned right after
# "RestrictionAndLigationSource" and right before "CRISPRSource",
# mirroring the other *Source sheets (headers in row 1 + a dropdown
# data validation on the first column).
# See: https://github.com/manulera/ShareYourCloning_backend/issues/215

$wb = $excel.ActiveWorkbook

$new = $wb.Worksheets.Add()
$new.Name = "GatewaySource"

# Re-fetch a live reference to the anchor sheet *after* Add(), then move
# the new sheet to sit right after it (i.e. right before CRISPRSource).
$afterSheet = $wb.Worksheets.Item("RestrictionAndLigationSource")
$new.Move($null, $afterSheet)

# Moving the sheet re-seats the original $new handle to whatever now
# occupies its old slot, so re-resolve a fresh reference by name before
# touching it any further.
$ws = $wb.Worksheets.Item("GatewaySource")

# Match the outline/page-setup look-and-feel of its sibling *Source sheets.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row, matching the column layout used by the other source sheets.
$ws.Range("A1").Value = "reaction_type"
$ws.Range("B1").Value = "circular"
$ws.Range("C1").Value = "assembly"
$ws.Range("D1").Value = "input"
$ws.Range("E1").Value = "output"
$ws.Range("F1").Value = "type"
$ws.Range("G1").Value = "output_name"
$ws.Range("H1").Value = "id"

# Dropdown validation on column A (reaction_type): LR or BR.
$rng = $ws.Range("A2:A1048576")
$rng.Validation.Add(3, 1, 1, """LR,BR""")
